# Insert a new data row at row 529 (shifting existing rows 529:600 down to
# 530:601) and populate it with the new record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 529..600 down to 530..601, leaving row 529 free for new data.
$ws.Rows.Item(529).Insert()

# Populate the newly inserted row 529 with its data.
$ws.Range("A529").Value = 9
$ws.Range("B529").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C529").Value = "Metropolitana"
$ws.Range("D529").Value = 45131
$ws.Range("E529").Value = 13
$ws.Range("F529").Value = 100112032
$ws.Range("G529").Value = "Zapallo italiano"
$ws.Range("H529").Value = "Sin especificar"
$ws.Range("I529").Value = "Primera"
$ws.Range("J529").Value = 70
$ws.Range("K529").Value = 17000
$ws.Range("L529").Value = 18000
$ws.Range("M529").Value = 17500
$ws.Range("N529").Value = "`$/caja 60 unidades"
$ws.Range("O529").Value = "Región de Arica y Parinacota"
$ws.Range("P529").Value = 292
$ws.Range("Q529").Value = 60
$ws.Range("R529").Value = "Hortaliza"
